$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("L40").ClearContents()
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = 0
$ws.Range("H43").Value = 51873.227
$ws.Range("I43").Value = 2090.125
$ws.Range("J43").Value = 80320.71000000001
$ws.Range("K43").Value = 2090.125
$ws.Range("L43").Value = 80320.71000000001
$ws.Range("M43").Value = -2021.125
$ws.Range("N43").Value = -80458.71000000001
$ws.Range("H76").Value = 7869.467
$ws.Range("I76").Value = 11583.385
$ws.Range("J76").Value = 5029.4116
$ws.Range("K76").Value = 11583.385
$ws.Range("L76").Value = 5029.4116
$ws.Range("M76").Value = -11268.385
$ws.Range("N76").Value = -5659.4116
$ws.Range("H79").Value = 7869.467
$ws.Range("I79").Value = 11583.385
$ws.Range("J79").Value = 5029.4116
$ws.Range("K79").Value = 11583.385
$ws.Range("L79").Value = 5029.4116
$ws.Range("M79").Value = -10491.385
$ws.Range("N79").Value = -7213.4116
$ws.Range("H113").Value = 4118.6294
$ws.Range("I113").Value = 3833.75
$ws.Range("J113").Value = 4346.533
$ws.Range("K113").Value = 3833.75
$ws.Range("L113").Value = 4346.533
$ws.Range("M113").Value = -579.75
$ws.Range("N113").Value = -10854.533
$ws.Range("H116").Value = 108799.25
$ws.Range("I116").Value = 143099
$ws.Range("K116").Value = 143099
$ws.Range("M116").Value = -139657
$ws.Range("H121").Value = 1176.0834
$ws.Range("J121").Value = 1232.091
$ws.Range("L121").Value = 3696.273
$ws.Range("N121").Value = -7190.272999999999
$ws.Range("H132").Value = 6948.478
$ws.Range("I132").Value = 3489.7222
$ws.Range("J132").Value = 19400
$ws.Range("K132").Value = 10469.1666
$ws.Range("L132").Value = 58200
$ws.Range("M132").Value = -7939.1666
$ws.Range("N132").Value = -63260

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1152.0869
$ws.Range("I45").Value = 1000.61536
$ws.Range("J45").Value = 1349
$ws.Range("K45").Value = 1000.61536
$ws.Range("L45").Value = 1349
$ws.Range("M45").Value = -623.61536
$ws.Range("N45").Value = -2103
$ws.Range("H74").Value = 2001.5957
$ws.Range("I74").Value = 1236.875
$ws.Range("J74").Value = 3633
$ws.Range("K74").Value = 1236.875
$ws.Range("L74").Value = 3633
$ws.Range("M74").Value = -362.875
$ws.Range("N74").Value = -5381
$ws.Range("H77").Value = 2001.5957
$ws.Range("I77").Value = 1236.875
$ws.Range("J77").Value = 3633
$ws.Range("K77").Value = 6184.375
$ws.Range("L77").Value = 18165
$ws.Range("M77").Value = -1816.375
$ws.Range("N77").Value = -26901
$ws.Range("H88").Value = 3399.5557
$ws.Range("I88").Value = 3685.1428
$ws.Range("J88").Value = 2400
$ws.Range("K88").Value = 3685.1428
$ws.Range("L88").Value = 2400
$ws.Range("M88").Value = -3279.1428
$ws.Range("N88").Value = -3212
$ws.Range("H91").Value = 3399.5557
$ws.Range("I91").Value = 3685.1428
$ws.Range("J91").Value = 2400
$ws.Range("K91").Value = 3685.1428
$ws.Range("L91").Value = 2400
$ws.Range("M91").Value = -2281.1428
$ws.Range("N91").Value = -5208

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2391.3816
$ws.Range("I105").Value = 2195.8708
$ws.Range("J105").Value = 3257.2144
$ws.Range("K105").Value = 2195.8708
$ws.Range("L105").Value = 3257.2144
$ws.Range("M105").Value = -448.8708000000001
$ws.Range("N105").Value = -6751.2144

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H18").Value = 37842
$ws.Range("I18").Value = 25000
$ws.Range("J18").Value = 50684
$ws.Range("K18").Value = 25000
$ws.Range("L18").Value = 50684
$ws.Range("M18").Value = -24770
$ws.Range("N18").Value = -51144
$ws.Range("H43").Value = 30552.334
$ws.Range("J43").Value = 30552.334
$ws.Range("L43").Value = 30552.334
$ws.Range("N43").Value = -30920.334
$ws.Range("H101").Value = 30552.334
$ws.Range("J101").Value = 30552.334
$ws.Range("L101").Value = 30552.334
$ws.Range("N101").Value = -37042.334
$ws.Range("H107").Value = 443.57144
$ws.Range("I107").Value = 335
$ws.Range("K107").Value = 335
$ws.Range("M107").Value = 1585
$ws.Range("H109").Value = 29995
$ws.Range("J109").Value = 29995
$ws.Range("L109").Value = 29995
$ws.Range("N109").Value = -32075

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 547.3103599999999
$ws.Range("I5").Value = 398.83334
$ws.Range("J5").Value = 1260
$ws.Range("K5").Value = 1196.50002
$ws.Range("L5").Value = 3780
$ws.Range("M5").Value = -1084.50002
$ws.Range("N5").Value = -4004
$ws.Range("H131").Value = 2000.0385
$ws.Range("I131").Value = 2126.9
$ws.Range("J131").Value = 1920.75
$ws.Range("K131").Value = 6380.700000000001
$ws.Range("L131").Value = 5762.25
$ws.Range("M131").Value = -1340.700000000001
$ws.Range("N131").Value = -15842.25
$ws.Range("H135").Value = 547.3103599999999
$ws.Range("I135").Value = 398.83334
$ws.Range("J135").Value = 1260
$ws.Range("K135").Value = 3589.50006
$ws.Range("L135").Value = 11340
$ws.Range("M135").Value = -1054.50006
$ws.Range("N135").Value = -16410
$ws.Range("H140").Value = 4057.5386
$ws.Range("I140").Value = 4194.8
$ws.Range("K140").Value = 12584.4
$ws.Range("M140").Value = -7404.400000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H69").Value = 20000
$ws.Range("J69").Value = 20000
$ws.Range("L69").Value = 20000
$ws.Range("N69").Value = -21498
$ws.Range("H72").Value = 20000
$ws.Range("J72").Value = 20000
$ws.Range("L72").Value = 60000
$ws.Range("N72").Value = -67488
$ws.Range("H101").Value = 34254.5
$ws.Range("J101").Value = 34254.5
$ws.Range("L101").Value = 34254.5
$ws.Range("N101").Value = -40744.5
$ws.Range("H107").Value = 7078.933
$ws.Range("I107").Value = 10257.8
$ws.Range("J107").Value = 721.2
$ws.Range("K107").Value = 10257.8
$ws.Range("L107").Value = 721.2
$ws.Range("M107").Value = -8337.799999999999
$ws.Range("N107").Value = -4561.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 372.8
$ws.Range("I22").Value = 317.66666
$ws.Range("J22").Value = 455.5
$ws.Range("K22").Value = 317.66666
$ws.Range("L22").Value = 455.5
$ws.Range("M22").Value = -22.66665999999998
$ws.Range("N22").Value = -1045.5
$ws.Range("H27").Value = 372.8
$ws.Range("I27").Value = 317.66666
$ws.Range("J27").Value = 455.5
$ws.Range("K27").Value = 317.66666
$ws.Range("L27").Value = 455.5
$ws.Range("M27").Value = -210.66666
$ws.Range("N27").Value = -669.5
$ws.Range("H63").Value = 25999
$ws.Range("J63").Value = 25999
$ws.Range("L63").Value = 25999
$ws.Range("N63").Value = -27497
$ws.Range("H66").Value = 25999
$ws.Range("J66").Value = 25999
$ws.Range("L66").Value = 77997
$ws.Range("N66").Value = -85485
$ws.Range("H68").Value = 66668416
$ws.Range("I68").Value = 76924370
$ws.Range("J68").Value = 4700
$ws.Range("K68").Value = 76924370
$ws.Range("L68").Value = 4700
$ws.Range("M68").Value = -76923621
$ws.Range("N68").Value = -6198
$ws.Range("H71").Value = 66668416
$ws.Range("I71").Value = 76924370
$ws.Range("J71").Value = 4700
$ws.Range("K71").Value = 384621850
$ws.Range("L71").Value = 23500
$ws.Range("M71").Value = -384618106
$ws.Range("N71").Value = -30988
$ws.Range("H122").Value = 3056.88
$ws.Range("I122").Value = 2917
$ws.Range("J122").Value = 3499.8333
$ws.Range("K122").Value = 8751
$ws.Range("L122").Value = 10499.4999
$ws.Range("M122").Value = -6301
$ws.Range("N122").Value = -15399.4999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 3339
$ws.Range("I45").Value = 1569
$ws.Range("J45").Value = 3781.5
$ws.Range("K45").Value = 1569
$ws.Range("L45").Value = 3781.5
$ws.Range("M45").Value = -1078
$ws.Range("N45").Value = -4763.5
$ws.Range("H62").Value = 3650.6
$ws.Range("I62").Value = 2556.889
$ws.Range("J62").Value = 4545.4546
$ws.Range("K62").Value = 2556.889
$ws.Range("L62").Value = 4545.4546
$ws.Range("M62").Value = -1932.889
$ws.Range("N62").Value = -5793.4546
$ws.Range("H65").Value = 3650.6
$ws.Range("I65").Value = 2556.889
$ws.Range("J65").Value = 4545.4546
$ws.Range("K65").Value = 12784.445
$ws.Range("L65").Value = 22727.273
$ws.Range("M65").Value = -9664.445
$ws.Range("N65").Value = -28967.273
